$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.713.64"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "2.331.24"
$ws.Range("E3").Value = "  +4.81%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "270.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "94.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.28%  "
$ws.Range("E7").Value = "  +1.03%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.619"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.59%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.01"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.88%  "
$ws.Range("E13").Value = "  +0.36%  "
$ws.Range("D14").Value = "2.681.47"
$ws.Range("E14").Value = "  +4.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.857"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.26%  "
$ws.Range("D17").Value = "2.347.09"
$ws.Range("E17").Value = "  +4.97%  "
$ws.Range("D18").Value = "43.650.91"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000107"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.95%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.66%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.30%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "238.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +9.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.53"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.34%  "
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.45"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.58"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.56"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "172.02"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0894"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("E35").Value = "  +2.87%  "
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("E37").Value = "  -3.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.35"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.235"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +14.00%  "
$ws.Range("E42").Value = "  +22.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.01"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.03%  "
$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "61.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.42%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "9.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.103"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "100.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("D50").Value = "2.556.54"
$ws.Range("E50").Value = "  +4.73%  "
$ws.Range("B51").Value = "TheGraph"
$ws.Range("C51").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.182"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.27%  "
